$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.057577835014504
$ws.Range("D2").Value = 1.055673677659529
$ws.Range("E2").Value = 1.063512524052167
$ws.Range("F2").Value = 1.073908068995313
$ws.Range("I2").Value = 1.050135973658008
$ws.Range("J2").Value = 1.06257353045938
$ws.Range("K2").Value = 1.058413488004615
$ws.Range("L2").Value = 1.066230962235443
$ws.Range("M2").Value = 1.076598684660691
$ws.Range("N2").Value = 1.064082506751519

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.05884592847386
$ws.Range("D3").Value = 1.056638677713175
$ws.Range("E3").Value = 1.064718401152645
$ws.Range("F3").Value = 1.075363787123258
$ws.Range("I3").Value = 1.050566234037704
$ws.Range("J3").Value = 1.063492617072769
$ws.Range("K3").Value = 1.059191589666079
$ws.Range("L3").Value = 1.067250884005787
$ws.Range("M3").Value = 1.077869858144718
$ws.Range("N3").Value = 1.065002898573321

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.059665679180933
$ws.Range("D4").Value = 1.057262313641685
$ws.Range("E4").Value = 1.065498174870204
$ws.Range("F4").Value = 1.07630565280832
$ws.Range("I4").Value = 1.050842891303345
$ws.Range("J4").Value = 1.064086020983154
$ws.Range("K4").Value = 1.059693671166548
$ws.Range("L4").Value = 1.067909761344243
$ws.Range("M4").Value = 1.078691787287593
$ws.Range("N4").Value = 1.065597145185323

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.060010116634057
$ws.Range("D5").Value = 1.057524304545659
$ws.Range("E5").Value = 1.065825872755247
$ws.Range("F5").Value = 1.076701597633089
$ws.Range("I5").Value = 1.050958780366423
$ws.Range("J5").Value = 1.064335177618466
$ws.Range("K5").Value = 1.05990441248193
$ws.Range("L5").Value = 1.068186497495338
$ws.Range("M5").Value = 1.079037184923474
$ws.Range("N5").Value = 1.065846655651635

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.06006793840795
$ws.Range("D6").Value = 1.057568283142104
$ws.Range("E6").Value = 1.065880887780758
$ws.Range("F6").Value = 1.076768077716617
$ws.Range("I6").Value = 1.050978214196324
$ws.Range("J6").Value = 1.064376993984209
$ws.Range("K6").Value = 1.059939777362015
$ws.Range("L6").Value = 1.068232947787014
$ws.Range("M6").Value = 1.079095170517222
$ws.Range("N6").Value = 1.065888531401415

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.059670282296508
$ws.Range("D7").Value = 1.057265815103971
$ws.Range("E7").Value = 1.065502554049223
$ws.Range("F7").Value = 1.076310943498463
$ws.Range("I7").Value = 1.050844441458158
$ws.Range("J7").Value = 1.064089351442828
$ws.Range("K7").Value = 1.059696488411278
$ws.Range("L7").Value = 1.067913460108284
$ws.Range("M7").Value = 1.078696403058877
$ws.Range("N7").Value = 1.065600480374631

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.058006558945654
$ws.Range("D8").Value = 1.055999967386921
$ws.Range("E8").Value = 1.063920163126706
$ws.Range("F8").Value = 1.074400054963134
$ws.Range("I8").Value = 1.050281745291878
$ws.Range("J8").Value = 1.062884412086414
$ws.Range("K8").Value = 1.058676742647052
$ws.Range("L8").Value = 1.066575874891593
$ws.Range("M8").Value = 1.077028411082802
$ws.Range("N8").Value = 1.064393829866123

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.05506865325879
$ws.Range("D9").Value = 1.053763291752551
$ws.Range("E9").Value = 1.061127753524023
$ws.Range("F9").Value = 1.071032025508133
$ws.Range("I9").Value = 1.049276741753935
$ws.Range("J9").Value = 1.060751038804106
$ws.Range("K9").Value = 1.056868991787127
$ws.Range("L9").Value = 1.064210482763742
$ws.Range("M9").Value = 1.074084405305178
$ws.Range("N9").Value = 1.062257426949053

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.053105649747332
$ws.Range("D10").Value = 1.05226795375437
$ws.Range("E10").Value = 1.059263246111482
$ws.Range("F10").Value = 1.068785885287661
$ws.Range("I10").Value = 1.048597606112534
$ws.Range("J10").Value = 1.059321833599895
$ws.Range("K10").Value = 1.055656418303762
$ws.Range("L10").Value = 1.062627736229412
$ws.Range("M10").Value = 1.072118293122863
$ws.Range("N10").Value = 1.06082619210952

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.052254553809115
$ws.Range("D11").Value = 1.051619425921153
$ws.Range("E11").Value = 1.058455163170101
$ws.Range("F11").Value = 1.067813038937472
$ws.Range("I11").Value = 1.048301347720986
$ws.Range("J11").Value = 1.058701286947778
$ws.Range("K11").Value = 1.055129575346969
$ws.Range("L11").Value = 1.061940970272391
$ws.Range("M11").Value = 1.071266078855766
$ws.Range("N11").Value = 1.060204764209979

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.051938248835412
$ws.Range("D12").Value = 1.051378375490016
$ws.Range("E12").Value = 1.058154890414615
$ws.Range("F12").Value = 1.06745163758984
$ws.Range("I12").Value = 1.048190973719089
$ws.Range("J12").Value = 1.058470531071299
$ws.Range("K12").Value = 1.054933610841614
$ws.Range("L12").Value = 1.061685657075845
$ws.Range("M12").Value = 1.070949392751886
$ws.Range("N12").Value = 1.059973680633687

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.052006105092086
$ws.Range("D13").Value = 1.051430088831571
$ws.Range("E13").Value = 1.058219305194097
$ws.Range("F13").Value = 1.067529161426044
$ws.Range("I13").Value = 1.048214664297787
$ws.Range("J13").Value = 1.058520040679994
$ws.Range("K13").Value = 1.054975658223603
$ws.Range("L13").Value = 1.061740432532218
$ws.Range("M13").Value = 1.071017329210795
$ws.Range("N13").Value = 1.060023260551705

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.052228411435833
$ws.Range("D14").Value = 1.051599503861936
$ws.Range("E14").Value = 1.058430344893545
$ws.Range("F14").Value = 1.067783166292928
$ws.Range("I14").Value = 1.04829223092236
$ws.Range("J14").Value = 1.058682217860306
$ws.Range("K14").Value = 1.055113382414643
$ws.Range("L14").Value = 1.061919870465867
$ws.Range("M14").Value = 1.071239904278811
$ws.Range("N14").Value = 1.060185668042215

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.052365359045951
$ws.Range("D15").Value = 1.051703864997432
$ws.Range("E15").Value = 1.058560358096628
$ws.Range("F15").Value = 1.067939661242965
$ws.Range("I15").Value = 1.048339978441834
$ws.Range("J15").Value = 1.058782106396133
$ws.Range("K15").Value = 1.055198202771559
$ws.Range("L15").Value = 1.062030399142254
$ws.Range("M15").Value = 1.071377022001067
$ws.Range("N15").Value = 1.060285698431221

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.053162110577743
$ws.Range("D16").Value = 1.052310972350795
$ws.Range("E16").Value = 1.05931686002642
$ws.Range("F16").Value = 1.068850444150411
$ws.Range("I16").Value = 1.048617221559369
$ws.Range("J16").Value = 1.059362981337101
$ws.Range("K16").Value = 1.055691345191632
$ws.Range("L16").Value = 1.062673284279302
$ws.Range("M16").Value = 1.072174832975719
$ws.Range("N16").Value = 1.060867398281232

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.053661593492029
$ws.Range("D17").Value = 1.052691515407323
$ws.Range("E17").Value = 1.059791193222043
$ws.Range("F17").Value = 1.069421683056584
$ws.Range("I17").Value = 1.048790541810542
$ws.Range("J17").Value = 1.059726893821933
$ws.Range("K17").Value = 1.056000199201982
$ws.Range("L17").Value = 1.063076164850155
$ws.Range("M17").Value = 1.072675041459934
$ws.Range("N17").Value = 1.061231827563534

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.05395282735374
$ws.Range("D18").Value = 1.052913379891855
$ws.Range("E18").Value = 1.060067792709698
$ws.Range("F18").Value = 1.069754852987257
$ws.Range("I18").Value = 1.048891425568811
$ws.Range("J18").Value = 1.059938994836794
$ws.Range("K18").Value = 1.056180175924116
$ws.Range("L18").Value = 1.063311020892727
$ws.Range("M18").Value = 1.07296672050498
$ws.Range("N18").Value = 1.061444229786164

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.054052112715684
$ws.Range("D19").Value = 1.052989013090033
$ws.Range("E19").Value = 1.060162094050303
$ws.Range("F19").Value = 1.069868451390918
$ws.Range("I19").Value = 1.048925788582987
$ws.Range("J19").Value = 1.060011288250878
$ws.Range("K19").Value = 1.05624151413659
$ws.Range("L19").Value = 1.063391077536432
$ws.Range("M19").Value = 1.073066161388039
$ws.Range("N19").Value = 1.061516625865188

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.053608014707212
$ws.Range("D20").Value = 1.0526506970856
$ws.Range("E20").Value = 1.059740309135681
$ws.Range("F20").Value = 1.069360397053699
$ws.Range("I20").Value = 1.048771968026719
$ws.Range("J20").Value = 1.059687866337786
$ws.Range("K20").Value = 1.055967079961512
$ws.Range("L20").Value = 1.063032953813615
$ws.Range("M20").Value = 1.072621382555425
$ws.Range("N20").Value = 1.061192744655883

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.052162952488127
$ws.Range("D21").Value = 1.051549619749199
$ws.Range("E21").Value = 1.058368202160578
$ws.Range("F21").Value = 1.067708369414284
$ws.Range("I21").Value = 1.04826939862026
$ws.Range("J21").Value = 1.058634467848459
$ws.Range("K21").Value = 1.055072833590377
$ws.Range("L21").Value = 1.061867036505402
$ws.Range("M21").Value = 1.071174365252155
$ws.Range("N21").Value = 1.060137850219875

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.051253398475659
$ws.Range("D22").Value = 1.050856411457514
$ws.Range("E22").Value = 1.057504838143726
$ws.Range("F22").Value = 1.066669423430088
$ws.Range("I22").Value = 1.047951500592705
$ws.Range("J22").Value = 1.05797066466137
$ws.Range("K22").Value = 1.054509012890919
$ws.Range("L22").Value = 1.06113271650595
$ws.Range("M22").Value = 1.070263780284908
$ws.Range("N22").Value = 1.059473104356121

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.051735665210536
$ws.Range("D23").Value = 1.051223982099868
$ws.Range("E23").Value = 1.057962588122016
$ws.Range("F23").Value = 1.067220213756902
$ws.Range("I23").Value = 1.048120206221179
$ws.Range("J23").Value = 1.058322701543407
$ws.Range("K23").Value = 1.054808054822178
$ws.Range("L23").Value = 1.06152211416747
$ws.Range("M23").Value = 1.070746574736064
$ws.Range("N23").Value = 1.059825641170909

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.053632224981341
$ws.Range("D24").Value = 1.052669141436622
$ws.Range("E24").Value = 1.059763301683389
$ws.Range("F24").Value = 1.069388089634486
$ws.Range("I24").Value = 1.048780361372056
$ws.Range("J24").Value = 1.059705501681413
$ws.Range("K24").Value = 1.055982045653885
$ws.Range("L24").Value = 1.06305247944536
$ws.Range("M24").Value = 1.072645628964309
$ws.Range("N24").Value = 1.06121040504372

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.055828931397115
$ws.Range("D25").Value = 1.054342260681566
$ws.Range("E25").Value = 1.061850155806005
$ws.Range("F25").Value = 1.071902862262384
$ws.Range("I25").Value = 1.04953816296155
$ws.Range("J25").Value = 1.06130378164517
$ws.Range("K25").Value = 1.057337634857797
$ws.Range("L25").Value = 1.064823006113241
$ws.Range("M25").Value = 1.074846091761801
$ws.Range("N25").Value = 1.062810954748351
